$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append (Rob_HPP_J extracted RoB data), columns:
# A=source_df, C=Primary, D=Study type, E=Q1, F=Q2, G=Q3, H=Q4, I=Q5, J=Q6, K=Q7, L=TotalStars, M=Q8, O=done_by
$newRows = @(
    @("Rob_HPP_J","Barrowclough(2015)","longitudinal prospective","1","1","1","0","2","1","1","7","0","Johanna"),
    @("Rob_HPP_J","Barrowclough(2013)","longitudinal prospective","1","1","1","0","2","1","1","7","0","Johanna"),
    @("Rob_HPP_J","Fond(2019)","longitudinal prospective","1","1","1","0","2","0","1","6","0","Johanna"),
    @("Rob_HPP_J","Foti(2010)","longitudinal prospective","1","1","1","0","2","1","1","8","1","Johanna"),
    @("Rob_HPP_J","Baeza(2009)","longitudinal prospective","1","1","0","1","2","1","1","8","1","Johanna"),
    @("Rob_HPP_J","Buchy(2015)","longitudinal prospective","1","1","1","1","1","1","1","7","0","Johanna"),
    @("Rob_HPP_J","Zammit(2011)","longitudinal prospective","1","1","0","1","2","1","1","8","1","Johanna"),
    @("Rob_HPP_J","Arseneault(2002)","longitudinal prospective","1","1","0","0","2","1","1","7","1","Johanna"),
    @("Rob_HPP_J","Bechtold(2016)","longitudinal prospective","1","1","0","0","2","0","1","6","1","Johanna"),
    @("Rob_HPP_J","Dragt(2011)","longitudinal prospective","1","1","1","1","2","1","1","9","1","Johanna")
)

$startRow = 131
$colLetters = @("A","C","D","E","F","G","H","I","J","K","L","M","O")

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $cellRef = "$($colLetters[$j])$r"
        $ws.Range($cellRef).Value = $rowData[$j]
    }
}
